# Graph.xlsx edit: add a straight-line-distance heuristic column (B) between
# the city-name column (A) and the neighbour-distance tuples (previously
# B:E, now shifted to C:F), and fix two city-name typos
# ("lasi" -> "Iasi", "Efo rie" -> "Eforie") everywhere they occur.
#
# NOTE: column widths/metadata are left completely untouched (no structural
# column insert is used) -- only cell contents are written -- since the
# target workbook's <cols> definitions are identical to the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fill column B with the heuristic (straight-line distance to Bucharest)
#    value for every city row, rows 1..20.
$heuristics = @(366, 329, 374, 380, 244, 241, 242, 160, 193, 253, 176, 10, 0, 77, 80, 199, 226, 234, 151, 161)
for ($i = 0; $i -lt $heuristics.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $heuristics[$i]
}

# 2) Move the existing neighbour-distance tuples one column to the right:
#    old B,C,D,E -> new C,D,E,F (written right-to-left per row so a cell's
#    old value is never clobbered before it has been copied onward).
$ws.Range("E1").Value = "(Sibiu, 140, 50, 11)"
$ws.Range("D1").Value = "(Zerind, 75, 50, 10.2)"
$ws.Range("C1").Value = "(Timisoara, 118, 50, 13.4)"
$ws.Range("D2").Value = "(Lugoj, 111, 60, 9)"
$ws.Range("C2").Value = "(Arad, 118, 50, 10)"
$ws.Range("D3").Value = "(Oradea, 71, 50, 10)"
$ws.Range("C3").Value = "(Arad, 75, 50, 7)"
$ws.Range("D4").Value = "(Sibiu, 151, 50, 10)"
$ws.Range("C4").Value = "(Zerind, 71, 50, 10)"
$ws.Range("D5").Value = "(Mehadia, 70, 50, 13.4)"
$ws.Range("C5").Value = "(Timisoara, 111, 50, 13.4)"
$ws.Range("D6").Value = "(Dobreta, 75, 50, 13.4)"
$ws.Range("C6").Value = "(Lugoj, 70, 50, 13.4)"
$ws.Range("D7").Value = "(Craiova, 120, 50, 13.4)"
$ws.Range("C7").Value = "(Mehadia, 75, 50, 13.4)"
$ws.Range("E8").Value = "(Pitesti, 138, 50, 13.4)"
$ws.Range("D8").Value = "(Rimnicu Vilcea, 146, 50, 13.4)"
$ws.Range("C8").Value = "(Dobreta, 120, 50, 13.4)"
$ws.Range("E9").Value = "(Pitesti, 97, 50, 13.4)"
$ws.Range("D9").Value = "(Craiova, 146, 50, 13.4)"
$ws.Range("C9").Value = "(Sibiu, 80, 50, 13.4)"
$ws.Range("F10").Value = "(Rimnicu Vilcea, 80, 50, 13.4)"
$ws.Range("E10").Value = "(Fagaras, 99, 50, 13.4)"
$ws.Range("D10").Value = "(Oradea, 151, 50, 13.4)"
$ws.Range("C10").Value = "(Arad, 140, 50, 13.4)"
$ws.Range("D11").Value = "(Bucharest, 211, 50, 13.4)"
$ws.Range("C11").Value = "(Sibiu, 99, 50, 13.4)"
$ws.Range("E12").Value = "(Bucharest, 101, 50, 13.4)"
$ws.Range("D12").Value = "(Rimnicu Vilcea, 97, 50, 13.4)"
$ws.Range("C12").Value = "(Craiova, 138, 50, 13.4)"
$ws.Range("F13").Value = "(Urziceni, 85, 50, 13.4)"
$ws.Range("E13").Value = "(Giurgiu, 90, 50, 13.4)"
$ws.Range("D13").Value = "(Pitesti, 101, 50, 13.4)"
$ws.Range("C13").Value = "(Fagaras, 211, 50, 13.4)"
$ws.Range("C14").Value = "(Bucharest, 90, 50, 13.4)"
$ws.Range("E15").Value = "(Hirsova, 98, 50, 13.4)"
$ws.Range("D15").Value = "(Vaslui, 142, 50, 13.4)"
$ws.Range("C15").Value = "(Bucharest, 85, 50, 13.4)"
$ws.Range("C16").Value = "(Urziceni, 142, 50, 13.4)"
$ws.Range("D17").Value = "(Neamt, 87, 50, 13.4)"
$ws.Range("C17").Value = "(Vaslui, 92, 50, 13.4)"
$ws.Range("C19").Value = "(Urziceni, 98, 50, 13.4)"
$ws.Range("C20").Value = "(Hirsova, 86, 50, 13.4)"

# 3) Fix the misspelled city names. "Efo rie" -> "Eforie" is corrected
#    first, then "lasi" -> "Iasi" (matching the order in which the shared
#    strings were authored), writing the neighbour-tuple cell before the
#    plain city-name cell in each case.
$ws.Range("D19").Value = "(Eforie, 86, 50, 13.4)"
$ws.Range("A20").Value = "Eforie"
$ws.Range("A17").Value = "Iasi"
$ws.Range("C18").Value = "(Iasi, 87, 50, 13.4)"
$ws.Range("D16").Value = "(Iasi, 92, 50, 13.4)"

# 4) Update the sheet view: scrolled one column to the right, with a new
#    active selection cell.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D25").Select()
